# Singapore Registry Config - feedback updates
# - "paid_up_capital" group row: rename the (C) node list and (G) field key
#   from the old "amount_with_currency" concept to "amount_under_paid_up_capital".
# - add a new trailing "Group" row (21) for total_equity_shares /
#   shareholdings_summary / number_of_shares_under_paid_up_capital.
# - update the window selection to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ("paid_up_capital" group) -------------------------------------
$ws.Range("C5").Value = "amount_under_paid_up_capital,currency"

# --- New row 21 ("total_equity_shares" group) -----------------------------
# Clone formatting from row 17 (same per-cell style pattern: A/G wrapped
# with style 3, B/C/D/E unstyled) then overwrite the values/height.
$ws.Range("A17:E17").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("G17").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Rows.Item(21).RowHeight = 30

$ws.Range("A21").Value = "total_equity_shares"
$ws.Range("B21").Value = "Group"
$ws.Range("D21").Value = "shareholdings_summary"
$ws.Range("C21").Value = "number_of_shares_under_paid_up_capital"
$ws.Range("E21").Value = "total_equity_shares"
$ws.Range("G21").Value = "number_of_shares_under_paid_up_capital"

# --- Row 5 (G) finished last, to match original authoring order -----------
$ws.Range("G5").Value = "amount_under_paid_up_capital"

# --- Window / selection state ---------------------------------------------
# NOTE: the target also scrolls the viewport so column C is the leftmost
# visible column (topLeftCell="C1"); this headless runtime's ActiveWindow
# does not persist ScrollColumn/TopLeftCell into the saved sheetView, so
# only the active-cell/selection half of the view state can be reproduced.
[void]$ws.Range("G5").Select()

Write-Output "edit complete"
